# Revert the previously-added "第596回" meeting entry.
# That entry lives in row 2 (right under the header row), so removing the
# entire row shifts every following meeting (第595回 ... 第570回) up by one,
# shrinking the used range from A1:F28 back down to A1:F27.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(2).Delete()
